$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column (C) for all data rows
# from serial date 45233 (2023-11-03) to 45243 (2023-11-13).
foreach ($row in 2..12) {
    $ws.Cells.Item($row, 3).Value = 45243
}
